$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Type" column header and per-row classification values.
$ws.Range("I1").Value = "Type"

$types = @(
    "Moving",      # row 2
    "Non-moving",  # row 3
    "Criminal",    # row 4
    "Criminal",    # row 5
    "Moving",      # row 6
    "Moving",      # row 7
    "Non-moving",  # row 8
    "Moving",      # row 9
    "Criminal",    # row 10
    "Moving",      # row 11
    "Moving",      # row 12
    "Criminal",    # row 13
    "Moving",      # row 14
    "Moving",      # row 15
    "Moving",      # row 16
    "Moving",      # row 17
    "Moving",      # row 18
    "Non-moving",  # row 19
    "Non-moving",  # row 20
    "Moving",      # row 21
    "Moving",      # row 22
    "Moving"       # row 23
)

for ($i = 0; $i -lt $types.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $types[$i]
}

# Row 10's FRA column (H) changes from "U" to "Y".
$ws.Range("H10").Value = "Y"

# Update the sheet view to match the post-edit scroll/selection state.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I24").Select()
